# Correct spelling mistakes in the Stanford_Speech sheet of expression.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stanford_Speech")

# "… had literally never a spreadsheet." -> "… had literally never met a spreadsheet."
$ws.Range("C72").Value = [char]0x2026 + " had literally never met a spreadsheet."

# "...could step in forme or for my brother..." -> "...could step in for me or for my brother..."
$ws.Range("C84").Value = "My colleagues were tremendous about stepping in for me at the halls of the United Nations, but nobody could step in for me or for my brother at the hospital."

# "...all our smrats, and all our soul." -> "...all our smarts, and all our soul."
$ws.Range("C86").Value = [char]0x2026 + " has a profound responsibility to try with all our skill, all our smarts, and all our soul."

# Leave the final selection on C87 to match the saved view state
$ws.Range("C87").Select()
